$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3 ("FOGKŐ ELTÁVOLÍTÁS" text block),
# shifting all subsequent rows down by one.
$ws.Rows.Item(3).Insert()

# Populate the new header row with the title/text column labels.
$ws.Range("A3").Value = "Cím"
$ws.Range("B3").Value = "Szöveg"

# Style the new row: bold font for both cells, with wrap text on the
# text column (matching the rest of column B).
$ws.Range("A3").Font.Bold = $true
$ws.Range("B3").Font.Bold = $true
$ws.Range("B3").WrapText = $true

# Match the selection left behind in the saved file.
$ws.Range("A3:B3").Select()
